$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.901.83'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '3.534.26'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.83'
$ws.Range("E5").Value = '  -2.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.60'
$ws.Range("E6").Value = '  +3.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  -1.35%  '
$ws.Range("E9").Value = '  -5.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.651'
$ws.Range("E10").Value = '  -2.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.71'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000302'
$ws.Range("E12").Value = '  -2.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.51'
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").Value = '4.098.06'
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '610.05'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.87'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.18'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '70.115.71'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = '3.541.73'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.994'
$ws.Range("E21").Value = '  -0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.01'
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.29'
$ws.Range("E23").Value = '  +2.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.49'
$ws.Range("E24").Value = '  -3.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.61'
$ws.Range("E25").Value = '  -2.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.11'
$ws.Range("E26").Value = '  +1.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.90'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.59'
$ws.Range("E28").Value = '  -5.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.55'
$ws.Range("E29").Value = '  -3.15%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("B31").Value = 'dogwifhat'
$ws.Range("C31").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.30'
$ws.Range("E31").Value = '  +14.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.63'
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("E33").Value = '  -2.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.24'
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").Value = '0.0₃0850'
$ws.Range("E35").Value = '  +8.32%  '
$ws.Range("D36").Value = '3.734.43'
$ws.Range("E36").Value = '  +5.35%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.04'
$ws.Range("E38").Value = '  -4.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.64'
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.392'
$ws.Range("E40").Value = '  -2.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.54'
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '492.72'
$ws.Range("E42").Value = '  -8.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.133'
$ws.Range("E43").Value = '  -4.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0456'
$ws.Range("E44").Value = '  -3.02%  '
$ws.Range("E45").Value = '  -2.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.83'
$ws.Range("E46").Value = '  -4.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.63'
$ws.Range("E49").Value = '  -4.45%  '
$ws.Range("E50").Value = '  +4.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.51'
$ws.Range("E51").Value = '  -2.46%  '
